$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.452.07"
$ws.Range("E2").Value = "  +2.84%  "
$ws.Range("D3").Value = "1.842.46"
$ws.Range("E3").Value = "  +2.07%  "
$ws.Range("E4").Value = "  +0.30%  "
$ws.Range("D5").Value = "'230.77"
$ws.Range("E5").Value = "  +2.87%  "
$ws.Range("D6").Value = "'0.610"
$ws.Range("E6").Value = "  +2.13%  "
$ws.Range("E7").Value = "  +0.22%  "
$ws.Range("D8").Value = "'43.45"
$ws.Range("E8").Value = "  +12.86%  "
$ws.Range("D9").Value = "'0.309"
$ws.Range("E9").Value = "  +7.68%  "
$ws.Range("D10").Value = "'0.0699"
$ws.Range("E10").Value = "  +4.67%  "
$ws.Range("E11").Value = "  +3.83%  "
$ws.Range("D12").Value = "2.107.80"
$ws.Range("E12").Value = "  +1.95%  "
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").Value = "'11.33"
$ws.Range("E13").Value = "  +2.30%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.833.92"
$ws.Range("E14").Value = "  +0.99%  "
$ws.Range("D15").Value = "'0.675"
$ws.Range("E15").Value = "  +7.88%  "
$ws.Range("D16").Value = "'4.69"
$ws.Range("E16").Value = "  +7.47%  "
$ws.Range("D17").Value = "35.400.75"
$ws.Range("E17").Value = "  +2.74%  "
$ws.Range("D18").Value = "'70.07"
$ws.Range("E18").Value = "  +2.91%  "
$ws.Range("D19").Value = "0.0₃0797"
$ws.Range("E19").Value = "  +3.95%  "
$ws.Range("D20").Value = "'244.87"
$ws.Range("E20").Value = "  +1.60%  "
$ws.Range("D21").Value = "'12.11"
$ws.Range("E21").Value = "  +9.31%  "
$ws.Range("D22").Value = "'4.66"
$ws.Range("E22").Value = "  +14.02%  "
$ws.Range("E23").Value = "  +0.24%  "
$ws.Range("D25").Value = "'168.98"
$ws.Range("E25").Value = "  -1.07%  "
$ws.Range("D26").Value = "'7.91"
$ws.Range("E26").Value = "  +3.17%  "
$ws.Range("E28").Value = "  +1.94%  "
$ws.Range("D29").Value = "'1.53"
$ws.Range("E29").Value = "  +24.95%  "
$ws.Range("E30").Value = "  +0.26%  "
$ws.Range("D31").Value = "3.378.93"
$ws.Range("E31").Value = "  +39.07%  "
$ws.Range("D32").Value = "'0.0547"
$ws.Range("E32").Value = "  +6.68%  "
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").Value = "'4.07"
$ws.Range("E33").Value = "  +6.71%  "
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "'3.93"
$ws.Range("E34").Value = "  +4.75%  "
$ws.Range("E35").Value = "  +2.26%  "
$ws.Range("D36").Value = "'94.93"
$ws.Range("E36").Value = "  +14.95%  "
$ws.Range("D37").Value = "'0.686"
$ws.Range("E37").Value = "  +7.50%  "
$ws.Range("D38").Value = "1.346.63"
$ws.Range("E38").Value = "  +2.31%  "
$ws.Range("D39").Value = "'1.08"
$ws.Range("E39").Value = "  +3.12%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").Value = "'2.43"
$ws.Range("E40").Value = "  +5.72%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").Value = "'0.0194"
$ws.Range("E41").Value = "  +3.85%  "
$ws.Range("B42").Value = "InjectiveProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D42").Value = "'15.19"
$ws.Range("E42").Value = "  +10.84%  "
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").Value = "'1.01"
$ws.Range("E43").Value = "  +7.25%  "
$ws.Range("E44").Value = "  +4.08%  "
$ws.Range("E45").Value = "  +1.02%  "
$ws.Range("E46").Value = "  +0.14%  "
$ws.Range("D47").Value = "'6.22"
$ws.Range("E47").Value = "  +7.82%  "
$ws.Range("E48").Value = "  +1.29%  "
$ws.Range("D49").Value = "2.008.71"
$ws.Range("E49").Value = "  +2.05%  "
$ws.Range("E50").Value = "  +0.27%  "
$ws.Range("D51").Value = "'103.14"
$ws.Range("E51").Value = "  +0.85%  "
